# Applies the "Updated cryptos list" data refresh described by the commit diff.
# For every changed cell we set the literal text that should appear in the sheet.
# Column D sometimes holds plain-looking numbers (e.g. "571.96") that Excel would
# otherwise auto-convert to a float (losing the trailing zero / exact text). We
# force those to stay text the same way a user would: a leading apostrophe.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '69.761.52'
$ws.Range("E2").Value = '  +0.10%  '

# Row 3
$ws.Range("D3").Value = '3.534.76'
$ws.Range("E3").Value = '  -0.77%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''571.96'
$ws.Range("E5").Value = '  -0.68%  '

# Row 6
$ws.Range("D6").Value = '''184.10'
$ws.Range("E6").Value = '  -3.30%  '

# Row 7
$ws.Range("D7").Value = '3.535.30'
$ws.Range("E7").Value = '  -0.59%  '

# Row 8
$ws.Range("D8").Value = '''0.616'
$ws.Range("E8").Value = '  -2.40%  '

# Row 9
$ws.Range("E9").Value = '  +0.10%  '

# Row 10
$ws.Range("D10").Value = '''0.182'
$ws.Range("E10").Value = '  +3.07%  '

# Row 11
$ws.Range("D11").Value = '''0.643'
$ws.Range("E11").Value = '  -2.37%  '

# Row 12
$ws.Range("D12").Value = '''54.00'
$ws.Range("E12").Value = '  -4.49%  '

# Row 13
$ws.Range("E13").Value = '  +0.03%  '

# Row 14
$ws.Range("D14").Value = '''9.50'
$ws.Range("E14").Value = '  -2.54%  '

# Row 15
$ws.Range("D15").Value = '4.101.56'
$ws.Range("E15").Value = '  -0.93%  '

# Row 16
$ws.Range("D16").Value = '''19.39'
$ws.Range("E16").Value = '  -3.17%  '

# Row 17 (WrappedBTC -> WrappedEther)
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.534.73'
$ws.Range("E17").Value = '  -1.10%  '

# Row 18 (WrappedEther -> WrappedBTC)
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '69.627.51'
$ws.Range("E18").Value = '  +0.05%  '

# Row 19
$ws.Range("D19").Value = '''12.39'
$ws.Range("E19").Value = '  -0.75%  '

# Row 20
$ws.Range("E20").Value = '  -1.21%  '

# Row 21
$ws.Range("E21").Value = '  -0.15%  '

# Row 22
$ws.Range("D22").Value = '''506.08'
$ws.Range("E22").Value = '  +5.93%  '

# Row 23
$ws.Range("D23").Value = '''19.83'
$ws.Range("E23").Value = '  +1.35%  '

# Row 24
$ws.Range("D24").Value = '''4.97'
$ws.Range("E24").Value = '  -2.14%  '

# Row 25
$ws.Range("D25").Value = '''4.37'
$ws.Range("E25").Value = '  +0.84%  '

# Row 26
$ws.Range("D26").Value = '''94.69'
$ws.Range("E26").Value = '  +6.56%  '

# Row 27
$ws.Range("D27").Value = '''11.26'
$ws.Range("E27").Value = '  +1.36%  '

# Row 28
$ws.Range("E28").Value = '  -4.41%  '

# Row 29
$ws.Range("E29").Value = '  -0.82%  '

# Row 30
$ws.Range("D30").Value = '''31.46'
$ws.Range("E30").Value = '  -1.80%  '

# Row 31
$ws.Range("D31").Value = '''7.51'
$ws.Range("E31").Value = '  -2.76%  '

# Row 32
$ws.Range("D32").Value = '''12.48'
$ws.Range("E32").Value = '  +3.04%  '

# Row 33
$ws.Range("D33").Value = '''65.42'
$ws.Range("E33").Value = '  -0.88%  '

# Row 34
$ws.Range("E34").Value = '  -5.03%  '

# Row 35
$ws.Range("D35").Value = '''572.36'
$ws.Range("E35").Value = '  -4.30%  '

# Row 36
$ws.Range("D36").Value = '''3.17'
$ws.Range("E36").Value = '  +8.21%  '

# Row 37 (Dai -> InjectiveProtocol)
$ws.Range("B37").Value = 'InjectiveProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D37").Value = '''37.92'
$ws.Range("E37").Value = '  -5.15%  '

# Row 38
$ws.Range("E38").Value = '  +0.75%  '

# Row 39 (InjectiveProtocol -> Dai)
$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").Value = '''0.999'
$ws.Range("E39").Value = '  +0.30%  '

# Row 40
$ws.Range("D40").Value = '0.0₃0775'
$ws.Range("E40").Value = '  -3.26%  '

# Row 41 (dogwifhat -> Stacks)
$ws.Range("B41").Value = 'Stacks'
$ws.Range("C41").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D41").Value = '''3.41'
$ws.Range("E41").Value = '  -2.72%  '

# Row 42 (Stacks -> dogwifhat)
$ws.Range("B42").Value = 'dogwifhat'
$ws.Range("C42").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D42").Value = '''3.13'
$ws.Range("E42").Value = '  -0.22%  '

# Row 43
$ws.Range("E43").Value = '  -7.29%  '

# Row 44
$ws.Range("D44").Value = '''3.58'
$ws.Range("E44").Value = '  +5.75%  '

# Row 45
$ws.Range("E45").Value = '  -4.12%  '

# Row 46
$ws.Range("D46").Value = '''0.0446'
$ws.Range("E46").Value = '  +0.04%  '

# Row 47
$ws.Range("D47").Value = '3.165.27'
$ws.Range("E47").Value = '  -2.57%  '

# Row 48
$ws.Range("E48").Value = '  -1.49%  '

# Row 49
$ws.Range("D49").Value = '''0.135'
$ws.Range("E49").Value = '  -1.68%  '

# Row 50 (LidoDAOToken -> OceanProtocol)
$ws.Range("B50").Value = 'OceanProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean'
$ws.Range("D50").Value = '''1.48'
$ws.Range("E50").Value = '  +26.19%  '

# Row 51
$ws.Range("D51").Value = '''0.999'
$ws.Range("E51").Value = '  +0.03%  '
